# Apply the "contingencies with rene fine" edit:
#  - Extend header row (row 1) with two new columns: P1=14, Q1=15 (same formatting as existing header cells)
#  - For data rows 2-25: swap values in columns I/K and M/O (I:1->2, K:2->1, M:1->2, O:2->1)
#  - For data rows 2-25: add two new columns P=2, Q=2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the existing header cells (N1:O1) onto the new header cells (P1:Q1),
# since the bold/border/centered style used there is a direct cell format, not a named style.
$ws.Range("N1:O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new) -> 2
    $ws.Cells.Item($r, 17).Value = 2   # Q (new) -> 2
}
